# Remove the "cost table" (tabelamento de custos) block that had been added
# below the Gantt schedule: the section title in B15 ("TABELA DE CUSTOS")
# and the cost breakdown table in B17:J29 (headers, rows, and totals).
# This reverts the sheet back to just the activities schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma de Atividades")

# Section title above the cost table.
$ws.Range("B15").ClearContents()

# Cost table: header rows, data rows, sub-totals and the two observation
# lines underneath it.
$ws.Range("B17:J29").ClearContents()
